$wb = $excel.ActiveWorkbook

# --- Employment Data sheet: update source B/C columns (and D for 2002-2011) ---
$ws = $wb.Worksheets.Item("Employment Data")

$data = @(
    @{ row = 8;  B = 40.700000762939503; C = 38.430824279785199 },
    @{ row = 9;  B = 41.200000762939503; C = 37.8655395507813 },
    @{ row = 10; B = 40.5;               C = 37.259838104247997 },
    @{ row = 11; B = 43.799999237060497; C = 36.636264801025398 },
    @{ row = 12; B = 45.799999237060497; C = 36.009674072265597 },
    @{ row = 13; B = 43.799999237060497; C = 35.381263732910199 },
    @{ row = 14; B = 43;                 C = 34.751255035400398 },
    @{ row = 15; B = 42.099998474121101; C = 34.130805969238303 },
    @{ row = 16; B = 42.200000762939503; C = 33.534633636474602 },
    @{ row = 17; B = 41.599998474121101; C = 32.975082397460902 },
    @{ row = 18; B = 42.5;               C = 32.458274841308601 },
    @{ row = 19; B = 41.299999237060497; C = 31.985078811645501; D = 45808736 },
    @{ row = 20; B = 40;                 C = 31.556421279907202; D = 46409243 },
    @{ row = 21; B = 39.799999237060497; C = 31.171117782592798; D = 47019452 },
    @{ row = 22; B = 41.200000762939503; C = 30.827831268310501; D = 47639556 },
    @{ row = 23; B = 42.400001525878899; C = 30.524616241455099; D = 48269753 },
    @{ row = 24; B = 41.799999237060497; C = 30.261022567748999; D = 48910248 },
    @{ row = 25; B = 42.400001525878899; C = 30.0376682281494;   D = 49561256 },
    @{ row = 26; B = 40.400001525878899; C = 29.855115890502901; D = 50222996 },
    @{ row = 27; B = 38.700000762939503; C = 29.711223602294901; D = 50895698 },
    @{ row = 28; B = 38.599998474121101; C = 29.6039218902588;   D = 51579599 }
)

foreach ($item in $data) {
    $r = $item.row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    if ($item.ContainsKey("D")) {
        $ws.Cells.Item($r, 4).Value = $item.D
    }
}

# Employment Data sheet becomes the active/selected sheet with D8:D28 selected
$ws.Activate()
$ws.Range("D8:D28").Select()

# --- ZAData sheet: selection changes (J2:J22 still, but tab no longer selected) ---
$zaData = $wb.Worksheets.Item("ZAData")
$zaData.Range("J2:J22").Select()

# --- Workbook view window settings ---
$excel.ActiveWindow.WindowState = -4143
$wb.Windows.Item(1).WindowState = -4143
